# Fix bugs: remove stray semicolons in D column utterance text for the
# "against positive" phrase rows, on the "Utterances" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Utterances")

$ws.Range("D13").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Pessoalmente esta é a minha decisão favorita.  </prosody></prosody></prosody>"
$ws.Range("D27").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Pessoalmente prefiro a outra decisão.  </prosody></prosody></prosody>"
$ws.Range("D37").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Definitely this will help people.  </prosody></prosody></prosody>"
$ws.Range("D39").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Wonderful. With this decision you will save the people.  </prosody></prosody></prosody>"
$ws.Range("D41").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Personally this is myfavourite decision.  </prosody></prosody></prosody>"
$ws.Range("D52").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Definitely I disagree with you.  </prosody></prosody></prosody>"
$ws.Range("D53").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Wrong. With this decision people will suffer.  </prosody></prosody></prosody>"
$ws.Range("D55").Value = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'><Gaze(person3)>  Personally I prefer other decision.  </prosody></prosody></prosody>"

# Update the sheet view/selection state to match the saved workbook
# (select the whole TEXT column, scrolled back to the top of the sheet)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D1:D1048576").Select()
